# Update the WESM exposure data table (rows 2-25) with the revised
# ACTUAL_ENERGY (B), CONTESTABLE_ENERGY (C), TOTAL_BCQ_NOMINATION (D)
# and WESM_EXPOSURE (E) figures. Row 19's WESM_EXPOSURE value is removed
# entirely (no longer computed/stored for that hour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 66539.3665
$ws.Range("C2").Value = 5516.6335
$ws.Range("D2").Value = 45000
$ws.Range("E2").Value = 16022.73300000001

$ws.Range("B3").Value = 63920.688
$ws.Range("C3").Value = 5387.312
$ws.Range("D3").Value = 45000
$ws.Range("E3").Value = 13533.376

$ws.Range("B4").Value = 60272.176
$ws.Range("C4").Value = 5339.824000000001
$ws.Range("D4").Value = 22500
$ws.Range("E4").Value = 32432.352

$ws.Range("B5").Value = 58962.181
$ws.Range("C5").Value = 5306.819000000001
$ws.Range("D5").Value = 22500
$ws.Range("E5").Value = 31155.36199999999

$ws.Range("B6").Value = 59797.17
$ws.Range("C6").Value = 5338.83
$ws.Range("D6").Value = 32500
$ws.Range("E6").Value = 21958.34

$ws.Range("B7").Value = 65642.2855
$ws.Range("C7").Value = 5422.714499999999
$ws.Range("D7").Value = 22500
$ws.Range("E7").Value = 37719.571

$ws.Range("B8").Value = 64173.893
$ws.Range("C8").Value = 5875.107
$ws.Range("D8").Value = 22500
$ws.Range("E8").Value = 35798.78599999999

$ws.Range("B9").Value = 73269.5725
$ws.Range("C9").Value = 7214.4275
$ws.Range("D9").Value = 22500
$ws.Range("E9").Value = 43555.14499999999

$ws.Range("B10").Value = 86257.982
$ws.Range("C10").Value = 8467.018
$ws.Range("D10").Value = 65000
$ws.Range("E10").Value = 12790.96400000001

$ws.Range("B11").Value = 90482.3835
$ws.Range("C11").Value = 13640.6165
$ws.Range("D11").Value = 65000
$ws.Range("E11").Value = 11841.76699999999

$ws.Range("B12").Value = 94548.98
$ws.Range("C12").Value = 15763.02
$ws.Range("D12").Value = 65000
$ws.Range("E12").Value = 13785.95999999999

$ws.Range("B13").Value = 96248.9605
$ws.Range("C13").Value = 15072.0395
$ws.Range("D13").Value = 65000
$ws.Range("E13").Value = 16176.921

$ws.Range("B14").Value = 93999.3325
$ws.Range("C14").Value = 15259.6675
$ws.Range("D14").Value = 65000
$ws.Range("E14").Value = 13739.66500000001

$ws.Range("B15").Value = 97613.256
$ws.Range("C15").Value = 15815.744
$ws.Range("D15").Value = 65000
$ws.Range("E15").Value = 16797.51199999999

$ws.Range("B16").Value = 98586.6725
$ws.Range("C16").Value = 15955.3275
$ws.Range("D16").Value = 65000
$ws.Range("E16").Value = 17631.345

$ws.Range("B17").Value = 83311.9265
$ws.Range("C17").Value = 16114.0735
$ws.Range("D17").Value = 65000
$ws.Range("E17").Value = 2197.853000000003

$ws.Range("B18").Value = 76092.611
$ws.Range("C18").Value = 16384.389
$ws.Range("D18").Value = 65000
$ws.Range("E18").Value = -5291.777999999991

$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 15932.637
$ws.Range("D19").Value = 65000
$ws.Range("E19").ClearContents()

$ws.Range("B20").Value = 81334.0935
$ws.Range("C20").Value = 15036.9065
$ws.Range("D20").Value = 65000
$ws.Range("E20").Value = 1297.187000000005

$ws.Range("B21").Value = 78578.4975
$ws.Range("C21").Value = 13330.5025
$ws.Range("D21").Value = 65000
$ws.Range("E21").Value = 247.9949999999953

$ws.Range("B22").Value = 77749.1885
$ws.Range("C22").Value = 11770.8115
$ws.Range("D22").Value = 65000
$ws.Range("E22").Value = 978.3770000000077

$ws.Range("B23").Value = 74837.9245
$ws.Range("C23").Value = 9312.075499999999
$ws.Range("D23").Value = 65000
$ws.Range("E23").Value = 525.8489999999947

$ws.Range("B24").Value = 54075.2015
$ws.Range("C24").Value = 6773.7985
$ws.Range("D24").Value = 65000
$ws.Range("E24").Value = -17698.59699999999

$ws.Range("B25").Value = 14561.85
$ws.Range("C25").Value = 5542.130999999999
$ws.Range("D25").Value = 65000
$ws.Range("E25").Value = -55980.281
